$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.751303672790527
$ws.Range("B1").Value = 2.738831520080566
$ws.Range("C1").Value = 3.394014596939087
$ws.Range("D1").Value = 1.299905776977539
$ws.Range("E1").Value = 0.8641474843025208
